# Apply the "Full positive set of addition tests" edit:
#  - Add the new AddingComputerHappyPathTest class name + 4 data-provider case
#    names as the automation-test linkage for TC-CD-01..04 on the
#    "Automation tests traceability" sheet.
#  - Widen columns D/E on that sheet to fit the new text.
#  - Restore the view/selection state that Excel would have saved after
#    the user's last interaction (leaving "Manual test cases" as the
#    selected/active sheet, scrolled near row 13).

$wb = $excel.ActiveWorkbook

$wsTrace = $wb.Worksheets.Item("Automation tests traceability")
$wsMain  = $wb.Worksheets.Item("Manual test cases")

$className = "\src\test\java\me\manzhos\tests\AddingComputerHappyPathTest"

# TC-CD-01 (row 3)
$wsTrace.Range("D3").Value = $className
$wsTrace.Range("E3").Value = "addComputerWithFieldsFromDataProviderTest (case 1 in dataprovider)"

# TC-CD-02 (row 4)
$wsTrace.Range("D4").Value = $className
$wsTrace.Range("E4").Value = "addComputerWithFieldsFromDataProviderTest (case 2 in dataprovider)"

# TC-CD-03 (row 5)
$wsTrace.Range("D5").Value = $className
$wsTrace.Range("E5").Value = "addComputerWithFieldsFromDataProviderTest (case 3 in dataprovider)"

# TC-CD-04 (row 6)
$wsTrace.Range("D6").Value = $className
$wsTrace.Range("E6").Value = "addComputerWithFieldsFromDataProviderTest (case 4 in dataprovider)"

# Widen columns D and E so the new, longer text fits (mirrors the author
# widening the columns after typing the new content).
$wsTrace.Columns.Item(4).ColumnWidth = 61.666666666666664
$wsTrace.Columns.Item(5).ColumnWidth = 63.49999999999999

# Update the saved selection/scroll state on each sheet, then leave the
# "Manual test cases" sheet active/selected last (it was tabSelected="1"
# both before and after the edit).
$wsTrace.Activate()
$wsTrace.Range("D26").Select()

$wsMain.Activate()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$wsMain.Range("D14").Select()
